$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style preparation via copy/paste-special (reuses existing style indices) ---
# Rows 23-32 take the A:s4 / B:s1 / C:s4 pattern already used by row 20
$ws.Range("A20:C20").Copy()
$ws.Range("A23:C32").PasteSpecial(-4122)

# Row 33 takes the A:s4 / B:(none) / C:s4 pattern already used by row 2
$ws.Range("A2:C2").Copy()
$ws.Range("A33:C33").PasteSpecial(-4122)

# Row 34 takes A:s4 (from row 2) but B/C revert to the default (no) style
$ws.Range("A2").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("B34:C34").Style = "Normal"

$excel.CutCopyMode = $false

# --- Cell values (rows 20-35) ---
$ws.Range("A20").Value = "Jonathan"
$ws.Range("B20").Value = "family_byers"
$ws.Range("C20").Value = "Me pidió que le cambiara el turno y acepté. Mamá, sabes que necesitamos el dinero"
$ws.Range("A21").Value = "Jonathan"
$ws.Range("B21").Value = "family_byers"
$ws.Range("C21").Value = "Confía en ti. No puedes dejar que te traten así"
$ws.Range("A22").Value = "Joyce Byers"
$ws.Range("B22").Value = "family_byers"
$ws.Range("C22").Value = "¿donde está tu hermano?¿No ha venido a dormir?"
$ws.Range("A23").Value = "Joyce Byers"
$ws.Range("B23").Value = "family_byers"
$ws.Range("C23").Value = "No puedes aceptar trabajos si yo tengo turno de noche, alguien tiene que encargarse de Will."
$ws.Range("A24").Value = "Joyce Byers"
$ws.Range("B24").Value = "hawkins"
$ws.Range("C24").Value = "Tranquila, no te preocupes, seguro que se habrá ido temprano al colegio."
$ws.Range("A25").Value = "Joyce Byers"
$ws.Range("B25").Value = "hawkins"
$ws.Range("C25").Value = "Tiene un par de amigos, pero la mayoría se ríen de el, de su ropa. Es un chico sensible, no es como los demás."
$ws.Range("A26").Value = "Joyce Byers"
$ws.Range("B26").Value = "hawkins"
$ws.Range("C26").Value = "Encuentra a mi hijo. "
$ws.Range("A27").Value = "Karen Wheeler"
$ws.Range("B27").Value = "hawkins"
$ws.Range("C27").Value = "Mañana hay cole, acabo de acostar a Holly, Terminad el fin de semana."
$ws.Range("A28").Value = "Nancy Wheeler"
$ws.Range("B28").Value = "friends_new"
$ws.Range("C28").Value = "Le gusto, pero no de esa forma. Nos enrollamos un par de veces."
$ws.Range("A29").Value = "Nancy Wheeler"
$ws.Range("B29").Value = "friends_new"
$ws.Range("C29").Value = "No puedo, tengo que estudiar, esos exámenes son imposibles."
$ws.Range("A30").Value = "Jim Hopper"
$ws.Range("B30").Value = "hawkins"
$ws.Range("C30").Value = "Estoy más guapo que tu mujer está mañana cuando la he dejado"
$ws.Range("A31").Value = "Dr. Owens"
$ws.Range("B31").Value = "hawkins"
$ws.Range("C31").Value = "Hemos cerrado esta zona siguiendo el protocolo de cuarentena"
$ws.Range("A32").Value = "Kali"
$ws.Range("B32").Value = "hawkins"
$ws.Range("C32").Value = "Me meteré en tu cabeza siempre que quiera"
$ws.Range("A33").Value = "Lucas"
$ws.Range("B33").Value = "friends"
$ws.Range("C33").Value = "Tu hermana ya era imbécil antes"
$ws.Range("A34").Value = "Lucas"
$ws.Range("B34").Value = "friends"
$ws.Range("C34").Value = "Lánzale una bola de fuego"
$ws.Range("A35").Value = "Will"
$ws.Range("B35").Value = "family_byers"
$ws.Range("C35").Value = "Es como esar en casa, pero es tan oscuro… es tan oscuro y vacío…"

# --- Clear stray chat_name (column B) values in rows 60-66, keep formatting ---
$ws.Range("B60").ClearContents()
$ws.Range("B61").ClearContents()
$ws.Range("B62").ClearContents()
$ws.Range("B63").ClearContents()
$ws.Range("B64").ClearContents()
$ws.Range("B65").ClearContents()
$ws.Range("B66").ClearContents()

# --- Update view state to mirror the author's last on-screen selection ---
$ws.Range("A266").Select()
